$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The handoff transform failed, so the status everywhere changes from
# "Ready for handoff" to "Handoff transform failed".
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"
$wsZhCn.Range("B2").Value = "Handoff transform failed"
$wsDeDe.Range("B2").Value = "Handoff transform failed"

function Reset-HandoffRow {
    param($ws)

    # Drop the hyperlink that pointed at the (now non-existent) handoff file
    # living in cell C2.
    $toDelete = @()
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$C$2') {
            $toDelete += $h
        }
    }
    foreach ($h in $toDelete) {
        $h.Delete()
    }

    # Fully remove the C2 cell (handoff file name) - no handoff happened.
    $ws.Range("C2").Clear()

    # Reset the handoff/handback datetimes and reason back to their
    # "nothing happened yet" defaults, matching row 3.
    $ws.Range("D2").Value = $ws.Range("D3").Value2
    $ws.Range("G2").Value = $ws.Range("G3").Value2
    $ws.Range("H2").Value = $ws.Range("H3").Value2
}

Reset-HandoffRow $wsZhCn
Reset-HandoffRow $wsDeDe
